# Clear the formula/value contents of C17:C24 while preserving their styles.
# These cells previously computed ages via "=2017-<birth year>" formulas;
# the calculation was removed/reworked so the cells should be blanked out.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C17:C24").ClearContents()

# Update the active selection to match the new selection C15:C30
$ws.Range("C15:C30").Select()
